$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.5808
$ws.Range("C21").Value = -12.3894
$ws.Range("C23").Value = -12.4166
$ws.Range("C25").Value = -13.0331
$ws.Range("D27").Value = -8.805900000000003
$ws.Range("D31").Value = -8.462400000000001
$ws.Range("D39").Value = -8.055600000000002
$ws.Range("D48").Value = -7.394099999999998
$ws.Range("D51").Value = -7.785999999999999
$ws.Range("D52").Value = -7.712799999999996
$ws.Range("C53").Value = -10.3509
$ws.Range("D55").Value = -8.251499999999997
$ws.Range("D56").Value = -8.0306
$ws.Range("C57").Value = -14.23199999999999
$ws.Range("D57").Value = -8.2514
$ws.Range("C59").Value = -12.7993
$ws.Range("C69").Value = -10.7612
$ws.Range("D73").Value = -7.662499999999999
$ws.Range("C79").Value = -10.62270000000001
$ws.Range("C83").Value = -13.97339999999999
$ws.Range("D89").Value = -5.994400000000004
$ws.Range("D90").Value = -8.091100000000003
$ws.Range("C93").Value = -11.20270000000001

$wb.Save()
